# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for the per-game rows to
# reflect strikeouts (K) instead of the previous Strike# metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 4
    10 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 3
    21 = 0
    22 = 3
    23 = 4
    24 = 1
    25 = 2
    26 = 6
    27 = 4
    28 = 3
    29 = 2
    30 = 1
    31 = 0
    32 = 0
    33 = 2
    35 = 1
    36 = 2
    37 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
